# Auto-generated edit script.
# Updates the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# columns for rows 2-51 on the active worksheet, matching the scheduled
# GitHub Actions data refresh described by the commit message:
#   "Updated cryptos list on Tue Jan  9 10:34:34 UTC 2024 with GitHub Actions"
#
# The source cells store these figures as literal text (inline/shared
# strings), e.g. "304.27" or "  +1.72%  ", not as numbers. Plain numeric
# looking strings (column D values without a second "." separator) are
# forced to Text format first so Excel's COM layer keeps them as the exact
# literal string instead of silently re-parsing them as a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.637.37"
$ws.Range("E2").Value = "  +6.00%  "
$ws.Range("D3").Value = "2.299.67"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.27"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.97"
$ws.Range("E6").Value = "  +11.54%  "
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +5.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.69"
$ws.Range("E10").Value = "  +10.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.39"
$ws.Range("E12").Value = "  +6.24%  "
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "2.650.01"
$ws.Range("E14").Value = "  +3.29%  "
$ws.Range("D15").Value = "2.300.69"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.86"
$ws.Range("E16").Value = "  +3.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.814"
$ws.Range("E17").Value = "  +4.73%  "
$ws.Range("D18").Value = "46.632.32"
$ws.Range("E18").Value = "  +6.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.09"
$ws.Range("E19").Value = "  +10.88%  "
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  +3.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.38"
$ws.Range("E22").Value = "  +2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.98"
$ws.Range("E23").Value = "  +5.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.90"
$ws.Range("E24").Value = "  +2.88%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +3.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.89"
$ws.Range("E27").Value = "  +10.69%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +4.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.02"
$ws.Range("E30").Value = "  +4.36%  "
$ws.Range("E31").Value = "  +12.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.67"
$ws.Range("E32").Value = "  +4.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "147.23"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.29"
$ws.Range("E35").Value = "  +15.56%  "
$ws.Range("E36").Value = "  +11.74%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("E38").Value = "  +5.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.05"
$ws.Range("E39").Value = "  +20.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.01"
$ws.Range("E40").Value = "  +10.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.37"
$ws.Range("E41").Value = "  +6.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0303"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  +9.77%  "
$ws.Range("D45").Value = "1.818.26"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.00"
$ws.Range("E46").Value = "  +20.91%  "
$ws.Range("E47").Value = "  +6.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.55"
$ws.Range("E48").Value = "  +9.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.92"
$ws.Range("E49").Value = "  +6.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.03"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.93"
$ws.Range("E51").Value = "  +5.44%  "
